$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.071.09"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "2.399.64"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.31"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.41"
$ws.Range("E6").Value = "  +4.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "2.412.66"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0968"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").Value = "  -4.74%  "
$ws.Range("D14").Value = "2.827.77"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").Value = "57.006.81"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.86"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "2.395.89"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.23"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "310.68"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("E22").Value = "  +4.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.86"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.23"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.377"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.42"
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.68"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("D31").Value = "0.0₃0724"
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.68"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.94"
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.995"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.84"
$ws.Range("E39").Value = "  +3.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.70"
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "133.42"
$ws.Range("E43").Value = "  +11.20%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.96"
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "252.93"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0489"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.05"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("E51").Value = "  +1.36%  "
